$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "In Translation" status text
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"

# Column width adjustments (target xml width 13.4101845877511 chars;
# COM ColumnWidth snaps to whole-pixel boundaries, so 12.5 is the closest
# achievable input, landing on the nearest representable width)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZh.Columns.Item(3).ColumnWidth = 12.5

$wsDe.Columns.Item(3).ColumnWidth = 12.5
